$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.869.75'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '2.408.29'
$ws.Range("E3").Value = '  +4.78%  '
$ws.Range("E4").Value = '  -0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '337.41'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +9.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.06'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -10.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.646'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.64%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.640'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.15'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -7.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0936'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.60'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.12%  '
$ws.Range("E13").Value = '  -3.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.98'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +8.95%  '
$ws.Range("E15").Value = '  +1.49%  '
$ws.Range("D16").Value = '2.766.81'
$ws.Range("E16").Value = '  +4.83%  '
$ws.Range("D17").Value = '2.403.76'
$ws.Range("E17").Value = '  +5.02%  '
$ws.Range("D18").Value = '42.975.91'
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.63'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +5.38%  '
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("E21").Value = '  +10.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '76.76'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '271.27'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +5.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.40'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.39'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +14.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.84'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.29'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +8.45%  '
$ws.Range("E29").Value = '  -1.63%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.10'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("E31").Value = '  -1.99%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0924'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.11%  '
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '36.44'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -7.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.03'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +4.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.135'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.81'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -5.68%  '
$ws.Range("E37").Value = '  -3.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.97'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -5.48%  '
$ws.Range("E39").Value = '  +3.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.91'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +12.24%  '
$ws.Range("E41").Value = '  +7.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.235'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.17'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.48%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.97'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +45.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '119.33'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +9.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.09'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.55'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.20'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.72%  '
$ws.Range("D50").Value = '1.650.72'
$ws.Range("E50").Value = '  +11.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.29'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.74%  '
